$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (Player, Position, Team) in the new row order.
$data = @(
    @('Aaron Wiggins',      'SG,SF',    'Oklahoma City Thunder'),
    @('Toumani Camara',     'SG,SF,PF', 'Portland Trail Blazers'),
    @('Devin Vassell',      'SG,SF',    'San Antonio Spurs'),
    @('Michael Porter Jr.', 'SF,PF',    'Denver Nuggets'),
    @('Jaden McDaniels',    'SF,PF',    'Minnesota Timberwolves'),
    @("De'Andre Hunter",    'SF,PF',    'Cleveland Cavaliers'),
    @('Josh Hart',          'SG,SF,PF', 'New York Knicks'),
    @('Domantas Sabonis',   'C',        'Sacramento Kings'),
    @('Alperen Sengün',     'C',        'Houston Rockets'),
    @('Guerschon Yabusele', 'PF,C',     'Philadelphia 76ers'),
    @('Kelly Oubre Jr.',    'SG,SF',    'Philadelphia 76ers'),
    @('Donovan Mitchell',   'PG,SG',    'Cleveland Cavaliers'),
    @('Malik Beasley',      'SG,SF',    'Detroit Pistons'),
    @('Dyson Daniels',      'PG,SG,SF', 'Atlanta Hawks'),
    @('Cam Thomas',         'SG,SF',    'Brooklyn Nets'),
    @('Julius Randle',      'PF,C',     'Minnesota Timberwolves'),
    @('Kristaps Porzingis', 'PF,C',     'Boston Celtics')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
